# Applies the "TIEMPO" -> "TEMPO" relabeling + tempo-code-to-text conversion
# described by the commit, on the "Respuestas de formulario 1" sheet, and
# un-hides the helper columns (B, C, G, H) that held the raw/intermediate
# data behind that computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Respuestas de formulario 1")

# Columns B:C and G:H were hidden helper columns; reveal them.
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(3).Hidden = $false
$ws.Columns.Item(7).Hidden = $false
$ws.Columns.Item(8).Hidden = $false

# Column I ("TIEMPO"/"TEMPO") held numeric tempo codes (1 = lento, 2 =
# rapido); replace them with the equivalent text labels row by row.
for ($r = 2; $r -le 61; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $val = $cell.Value2
    if ($val -eq 1) {
        $cell.Value = "Tempo lento"
    } elseif ($val -eq 2) {
        $cell.Value = "Tempo rapido"
    }
}

# The column header itself changes from "TIEMPO" to "TEMPO".
$ws.Range("I1").Value = "TEMPO"
